$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows whose option/description changed (row 12 & 13) ---
# Row 12 (A=9): Heart -> CurHeart, "최대 체력 {0} 증가" -> "체력 {0} 회복"
$ws.Range("B12").Value2 = "CurHeart"
$ws.Range("C12").Value2 = "체력 {0} 회복"

# Row 13 (A=10): DevilRoom/악마방.. -> MaxHeart / "최대 체력 {0} 증가"
$ws.Range("A13").Value2 = 10
$ws.Range("B13").Value2 = "MaxHeart"
$ws.Range("C13").Value2 = "최대 체력 {0} 증가"

# Row 14 (A=11): AngelRoom/천사방.. -> DevilRoom / "악마방 등장 확률 {0}% 증가"
$ws.Range("A14").Value2 = 11
$ws.Range("B14").Value2 = "DevilRoom"
$ws.Range("C14").Value2 = "악마방 등장 확률 {0}% 증가"

# Row 15 (A=12): BlackHeart/블랙하트.. -> AngelRoom / "천사방 등장 확률 {0}% 증가"
$ws.Range("A15").Value2 = 12
$ws.Range("B15").Value2 = "AngelRoom"
$ws.Range("C15").Value2 = "천사방 등장 확률 {0}% 증가"

# Row 16 (A=13): AttackPerCoin/코인당.. -> BlackHeart / "블랙 하트 {0} 증가"
$ws.Range("A16").Value2 = 13
$ws.Range("B16").Value2 = "BlackHeart"
$ws.Range("C16").Value2 = "블랙 하트 {0} 증가"

# Row 17 (A=14): SoulHeart/소울하트.. -> AttackPerCoin / "보유한 코인당 공격력 {0} 증가"
$ws.Range("A17").Value2 = 14
$ws.Range("B17").Value2 = "AttackPerCoin"
$ws.Range("C17").Value2 = "보유한 코인당 공격력 {0} 증가"

# Copy the formatting used by the B column (vertical-center + wrap) down to
# the newly-added rows (18-21) so they match the rest of the column.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B18:B21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New row 18 (A=15): SoulHeart / "소울 하트 {0} 증가"
$ws.Range("A18").Value2 = 15
$ws.Range("B18").Value2 = "SoulHeart"
$ws.Range("C18").Value2 = "소울 하트 {0} 증가"

# New row 19 (A=16): Coin / "코인 {0}개 획득"
$ws.Range("A19").Value2 = 16
$ws.Range("B19").Value2 = "Coin"
$ws.Range("C19").Value2 = "코인 {0}개 획득"

# New row 20 (A=17): Boom / "폭탄 {0}개 획득"
$ws.Range("A20").Value2 = 17
$ws.Range("B20").Value2 = "Boom"
$ws.Range("C20").Value2 = "폭탄 {0}개 획득"

# New row 21 (A=18): Key / "열쇠 {0}개 획득"
$ws.Range("A21").Value2 = 18
$ws.Range("B21").Value2 = "Key"
$ws.Range("C21").Value2 = "열쇠 {0}개 획득"

# --- Update the sheet view: zoom & selection ---
$ws.Select() | Out-Null
$ws.Range("C19").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 115
